# Update gh-pages to output generated at 456a3b4
#
# Two sheets ("展览" and "全部类型") each had their 2nd and 3rd data rows
# (old rows 3 & 4 - the "一周年超强巨制" full event and its "内场票" inner
# ticket variant) removed. Every later event row shifts up two positions,
# and the "想去人数" (interest count) column picks up refreshed numbers for
# a handful of rows. A couple of cells also flip between a numeric price
# and the literal "不可售" (not for sale) placeholder, but that value
# already lands correctly once the rows shift - only the "A" row index
# column and column F need to be re-applied after the shift.

$wb = $excel.ActiveWorkbook

# (sheet name, number of rows to delete at the top, F-column refresh list, last data row)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the two obsolete rows (old row 3 "...一周年超强巨制~" and old
    # row 4 "...内场票·赵成晨"). Deleting row 3 twice pulls every following
    # row up by two, which is exactly the shift the diff shows.
    $ws.Rows.Item(3).Delete()
    $ws.Rows.Item(3).Delete()

    # Work out how many data rows now remain (header row 1 + N data rows).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    if ($lastRow -lt 2) {
        $lastRow = $wb.Application.WorksheetFunction.Max(2, $lastRow)
    }

    # Column A held a plain sequential index (1, 2, 3, ...) - after the row
    # delete it still carries its *old* numbers, so renumber it to match.
    $idx = 1
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $idx
        $idx = $idx + 1
    }

    # Refresh the "想去人数" (F) counts for the rows whose number moved.
    $ws.Cells.Item(3, 6).Value = 994
    $ws.Cells.Item(4, 6).Value = 221
    $ws.Cells.Item(5, 6).Value = 1358
    $ws.Cells.Item(6, 6).Value = 8368
    $ws.Cells.Item(7, 6).Value = 49
    $ws.Cells.Item(9, 6).Value = 620
    $ws.Cells.Item(12, 6).Value = 3356
    $ws.Cells.Item(13, 6).Value = 42
    $ws.Cells.Item(15, 6).Value = 42
    $ws.Cells.Item(16, 6).Value = 859
    $ws.Cells.Item(21, 6).Value = 1906
}
